# SchedulerProgram/students.xlsx — "Student Details" sheet
#
# The live app logs each search/lookup by stamping the username into A2
# and the lookup timestamp into B2. The commit wires up real search for
# the ShoppingCart, so a fresh run of searches landed the username back
# to "kvw5270" and advanced the timestamp to the latest recorded run.
# The username column was also narrowed to better fit the (now shorter)
# search results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final recorded search result: username + timestamp of the last lookup.
$ws.Range("A2").Value = "kvw5270"
$ws.Range("B2").Value = "03/22/2020 23:09:36"

# Username column narrowed to fit the shorter search values.
$ws.Columns.Item(1).ColumnWidth = 9.33
